$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 3247
$ws.Cells.Item(7, 6).Value = 328
$ws.Cells.Item(8, 6).Value = 7596
$ws.Cells.Item(9, 6).Value = 80
$ws.Cells.Item(11, 6).Value = 5
$ws.Cells.Item(12, 6).Value = 19
$ws.Cells.Item(14, 6).Value = 645
$ws.Cells.Item(15, 6).Value = 1095
$ws.Cells.Item(16, 6).Value = 1030
$ws.Cells.Item(19, 6).Value = 1382
$ws.Cells.Item(21, 6).Value = 6010
$ws.Cells.Item(22, 6).Value = 20
$ws.Cells.Item(23, 6).Value = 2338
$ws.Cells.Item(24, 6).Value = 4148
$ws.Cells.Item(25, 6).Value = 2811
$ws.Cells.Item(26, 6).Value = 264
$ws.Cells.Item(29, 6).Value = 1017
$ws.Cells.Item(30, 6).Value = 261
$ws.Cells.Item(34, 6).Value = 1009
$ws.Cells.Item(35, 6).Value = 1008
$ws.Cells.Item(36, 6).Value = 76
$ws.Cells.Item(40, 6).Value = 187
$ws.Cells.Item(41, 6).Value = 17
$ws.Cells.Item(43, 6).Value = 363
$ws.Cells.Item(44, 6).Value = 297
$ws.Cells.Item(45, 6).Value = 1042
$ws.Cells.Item(48, 6).Value = 1739
$ws.Cells.Item(49, 6).Value = 55
$ws.Cells.Item(50, 6).Value = 311

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 592
$ws.Cells.Item(12, 6).Value = 69
$ws.Cells.Item(15, 6).Value = 168
$ws.Cells.Item(19, 6).Value = 142
$ws.Cells.Item(21, 6).Value = 32
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(27, 6).Value = 4369
$ws.Cells.Item(28, 6).Value = 4369
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(32, 6).Value = 48

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 23
$ws.Cells.Item(6, 6).Value = 1938
$ws.Cells.Item(8, 6).Value = 3010
$ws.Cells.Item(9, 6).Value = 1196
$ws.Cells.Item(10, 6).Value = 1241
$ws.Cells.Item(14, 6).Value = 8693
$ws.Cells.Item(15, 6).Value = 842

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 3247
$ws.Cells.Item(6, 6).Value = 1938
$ws.Cells.Item(7, 6).Value = 328
$ws.Cells.Item(8, 6).Value = 3010
$ws.Cells.Item(9, 6).Value = 7596
$ws.Cells.Item(10, 6).Value = 1196
$ws.Cells.Item(11, 6).Value = 1241
$ws.Cells.Item(12, 6).Value = 80
$ws.Cells.Item(15, 6).Value = 19
$ws.Cells.Item(17, 6).Value = 842
$ws.Cells.Item(19, 6).Value = 592
$ws.Cells.Item(20, 6).Value = 592
$ws.Cells.Item(21, 6).Value = 645
$ws.Cells.Item(22, 6).Value = 1095
$ws.Cells.Item(23, 6).Value = 1030
$ws.Cells.Item(24, 6).Value = 69
$ws.Cells.Item(26, 6).Value = 168
$ws.Cells.Item(27, 6).Value = 1382
$ws.Cells.Item(29, 6).Value = 6010
$ws.Cells.Item(30, 6).Value = 2338
$ws.Cells.Item(31, 6).Value = 4148
$ws.Cells.Item(32, 6).Value = 2811
$ws.Cells.Item(33, 6).Value = 264
$ws.Cells.Item(36, 6).Value = 261
$ws.Cells.Item(38, 6).Value = 76
$ws.Cells.Item(41, 6).Value = 187
$ws.Cells.Item(42, 6).Value = 142
$ws.Cells.Item(45, 6).Value = 297
$ws.Cells.Item(48, 6).Value = 1740
$ws.Cells.Item(49, 6).Value = 55
$ws.Cells.Item(50, 6).Value = 4369
$ws.Cells.Item(52, 6).Value = 48
